{"js": "// Bold the four key terms (\"importing\", \"exporting\", \"change\", \"Returning\")\n// that appear as spell-check-flagged words in the document body.\nconst body = context.document.body;\nconst terms = [\"importing\", \"exporting\", \"change\", \"Returning\"];\n\nfor (const term of terms) {\n  const results = body.search(term, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].font.bold = true;\n  }\n  await context.sync();\n}\n", "ps1": "# Bold the four key terms (\"importing\", \"exporting\", \"change\", \"Returning\")\n# that appear as spell-check-flagged words in the document body.\n$d = $word.ActiveDocument\n$terms = @(\"importing\", \"exporting\", \"change\", \"Returning\")\n\nforeach ($term in $terms) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($term, $true, $true)\n    if ($rng.Find.Found) {\n        $rng.Bold = 1\n    }\n}\n"}
